# ShellCheck fix (unused variable)
# Adds the new "Data" row (2023-10-16, row 83) that was captured by the
# varia/check_ssl_cert_stats.xlsx data-collection script, extends the
# "Data" Excel table accordingly, and tidies up the formatting that moves
# along with the new "latest row" (Stars/Forks General-number-format
# highlight) plus the column-width bump on the Total column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$lo = $ws.ListObjects.Item(1)

# --- 1. Grow the table by one row -----------------------------------------
# ListRows.Add() extends ListObject.Range (and therefore the table ref /
# autoFilter ref / worksheet dimension) from A1:AK82 to A1:AK83.
$null = $lo.ListRows.Add()

# --- 2. Inherit row 82's formatting into the new row 83 --------------------
# A straight range copy brings every number format (date, #,##0, the
# +/-#,##0 delta format, ...) down to row 83 before we overwrite the
# values; formulas are re-entered explicitly afterwards (see below) so the
# relative "-G82" / "-P82" style references stay correct for row 83.
$ws.Range("A82:AK82").Copy($ws.Range("A83:AK83"))

# --- 3. Write the new row's data -------------------------------------------
$ws.Range("A83").Value = 45257
$ws.Range("B83").Value = 347
$ws.Range("C83").Value = 129
$ws.Range("D83").Value = 114
$ws.Range("E83").Value = 282
$ws.Range("F83").Value = 233
$ws.Range("G83").Value = 5754
$ws.Range("H83").Formula = "=Data[[#This Row],[LoC]]-G82"
$ws.Range("I83").Value = 6920
$ws.Range("J83").Value = 1971
$ws.Range("K83").Value = 563
$ws.Range("L83").Value = 417
$ws.Range("M83").Value = 134
$ws.Range("N83").Value = 60
$ws.Range("O83").Value = 16
$ws.Range("P83").Formula = "=SUM(Data[[#This Row],[Shell]:[Bash]])"
$ws.Range("Q83").Formula = "=Data[[#This Row],[Total]]-P82"
$ws.Range("R83").Value = 2160
$ws.Range("S83").Value = 4583
$ws.Range("T83").Value = 71467
$ws.Range("U83").Value = 48982
$ws.Range("V83").Value = 2
$ws.Range("W83").Value = 1
$ws.Range("X83").Value = 285
$ws.Range("Y83").Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"
$ws.Range("Z83").Value = 0
$ws.Range("AA83").Value = 178
$ws.Range("AB83").Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"
$ws.Range("AC83").Value = 159
$ws.Range("AD83").Value = 166
$ws.Range("AE83").Value = 7
# Running / Failed / OK / Cancelled / GH runs were not recorded this time.
$ws.Range("AF83").ClearContents()
$ws.Range("AG83").ClearContents()
$ws.Range("AH83").ClearContents()
$ws.Range("AI83").ClearContents()
$ws.Range("AJ83").ClearContents()
$ws.Range("AK83").Formula = "=SUM(Data[[#This Row],[Running]:[GH runs]])"

# --- 4. Move the "latest row" Stars/Forks highlight down to row 83 ---------
# Row 82 carried an explicit (if visually-equivalent to plain) General
# number format on Stars/Forks; that marker now belongs to the new last
# row. Pull row 81's plain/unformatted Stars+Forks cells down over row 82
# to drop the marker, then re-set row 82's real values; and stamp row 83's
# Stars+Forks explicitly as General to carry the marker forward.
$ws.Range("B81:C81").Copy($ws.Range("B82:C82"))
$ws.Range("B82").Value = 346
$ws.Range("C82").Value = 129
$ws.Range("B83").NumberFormat = "General"
$ws.Range("C83").NumberFormat = "General"

# --- 5. Total column needed a little more room ------------------------------
$ws.Columns.Item(16).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# --- 6. Park the viewport/selection roughly where the author left it -------
$ws.Range("AB76").Select()
